# Add a new worksheet 'magapoke_2026-01-07' after the last existing sheet,
# populate it with the weekly ranking data, and style the header row
# the same way as the other ranking sheets (bold, thin box border,
# centered horizontally / top vertically).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = 'magapoke_2026-01-07'

# Match the page margins used by the other ranking sheets (0.75in / 1in / 0.5in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row
$ws.Cells.Item(1, 1).Value = 'rank'
$ws.Cells.Item(1, 2).Value = 'title'

$headerRange = $ws.Range('A1:B1')
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Ranking data rows (rank, title)
$data = @(
    @(1, 'イエティ、とある日々'),
    @(2, 'せいぶつ部の田辺くん'),
    @(3, '黄昏町プリズナーズ'),
    @(4, 'K-9~警視庁公安部公安第9課異能対策係~'),
    @(5, 'ハードワーカー中田'),
    @(6, 'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜'),
    @(7, '普通の本はありません！'),
    @(8, 'スルガメテオ'),
    @(9, 'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！'),
    @(10, 'ドリーム☆ジャンボ☆ガール'),
    @(11, 'アイドラトリィ'),
    @(12, '春くらり'),
    @(13, '黒月のイェルクナハト'),
    @(14, 'きゃわるり方程式'),
    @(15, 'ゼロとヒャク'),
    @(16, '屋根の下のアルテミス'),
    @(17, '篝家の８兄弟'),
    @(18, '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～'),
    @(19, 'MYS'),
    @(20, '歪みの虜'),
    @(21, '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～'),
    @(22, 'その青春'),
    @(23, '平成転生'),
    @(24, '君が監督！'),
    @(25, 'ともだちづくり'),
    @(26, '生きたがりの人狼'),
    @(27, '鳴るさんだぁ'),
    @(28, '明智ナンバーワン'),
    @(29, '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～'),
    @(30, 'JK Biker'),
    @(31, '夜鐘のキト'),
    @(32, '追放されなかった男　～二度目の人生は土下座から始まりました～'),
    @(33, 'ナキナギ'),
    @(34, 'じゅーくぼっくす'),
    @(35, '永久のユウグレ'),
    @(36, '〈小市民〉 春期限定いちごタルト事件'),
    @(37, 'GURU'),
    @(38, '卒業アルバムの彼女たち'),
    @(39, '花子狩り'),
    @(40, 'ハプスブルク家の華麗なる受難'),
    @(41, '白銀のキュイジーヌ～明治外交官の料理人～'),
    @(42, '人生逆転ダンジョン'),
    @(43, 'ナマイキ旭ちゃんをわからせたい'),
    @(44, '眠れる森のレガ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
}

